$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge existing merged header cells (row 1) so every cell can be written individually ---
$ws.Range("B1:P1").UnMerge()
$ws.Range("R1:AC1").UnMerge()
$ws.Range("AD1:AG1").UnMerge()
$ws.Range("AH1:AK1").UnMerge()
$ws.Range("AL1:AO1").UnMerge()
$ws.Range("AP1:AQ1").UnMerge()
$ws.Range("AR1:AS1").UnMerge()

# --- Clear old contents of the table (rows 1,2,4,5; row 3 stays blank as before) ---
$ws.Range("A1:AT2").ClearContents()
$ws.Range("A4:AT5").ClearContents()

# Row 1
$ws.Range('B1').Value2 = 'Biorefinery'
$ws.Range('R1').Value2 = 'Installed cost'
$ws.Range('S1').Value2 = 'Material cost'
$ws.Range('AK1').Value2 = 'Product sale'
$ws.Range('AO1').Value2 = 'Heating demand'
$ws.Range('AS1').Value2 = 'Cooling demand'
$ws.Range('AW1').Value2 = 'Power demand'
$ws.Range('AY1').Value2 = 'Utility cost'
$ws.Range('BA1').Value2 = 'TEA'
$ws.Range('BB1').Value2 = 'Biorefinery'

# Row 2
$ws.Range('B2').Value2 = 'Minimum selling price [$/kg]'
$ws.Range('C2').Value2 = 'Product yield [10^6 kg/yr]'
$ws.Range('D2').Value2 = 'Product purity [%]'
$ws.Range('E2').Value2 = 'Adjusted minimum selling price [$/kg]'
$ws.Range('F2').Value2 = 'Adjusted product yield [10^6 kg/yr]'
$ws.Range('G2').Value2 = 'Product recovery [%]'
$ws.Range('H2').Value2 = 'Total capital investment [10^6 $]'
$ws.Range('I2').Value2 = 'Annual operating cost [10^6 $/yr]'
$ws.Range('J2').Value2 = 'Annual material cost [10^6 $/yr]'
$ws.Range('K2').Value2 = 'Annual product sale [10^6 $/yr]'
$ws.Range('L2').Value2 = 'Annual electricity credit [10^6 $/yr]'
$ws.Range('M2').Value2 = 'Pretreatment - heating duty [kJ/kg]'
$ws.Range('N2').Value2 = 'Saccharification and Fermentation - heating duty [kJ/kg]'
$ws.Range('O2').Value2 = 'Separation - heating duty [kJ/kg]'
$ws.Range('P2').Value2 = 'Waste treatment - heating duty [kJ/kg]'
$ws.Range('Q2').Value2 = 'Product storage and pumping - heating duty [kJ/kg]'
$ws.Range('R2').Value2 = 'Check [10^6 $]'
$ws.Range('S2').Value2 = 'feedstock [10^6 $/yr]'
$ws.Range('T2').Value2 = 'sulfuric_acid_fresh [10^6 $/yr]'
$ws.Range('U2').Value2 = 'enzyme [10^6 $/yr]'
$ws.Range('V2').Value2 = 'boiler_chems [10^6 $/yr]'
$ws.Range('W2').Value2 = 'sulfuric_acid_fresh2 [10^6 $/yr]'
$ws.Range('X2').Value2 = 'baghouse_bag [10^6 $/yr]'
$ws.Range('Y2').Value2 = 'ammonia_fresh [10^6 $/yr]'
$ws.Range('Z2').Value2 = 'natural_gas [10^6 $/yr]'
$ws.Range('AA2').Value2 = 'CSL_fresh [10^6 $/yr]'
$ws.Range('AB2').Value2 = 'cooling_tower_chems [10^6 $/yr]'
$ws.Range('AC2').Value2 = 'lime_fresh [10^6 $/yr]'
$ws.Range('AD2').Value2 = 'hexanol_fresh [10^6 $/yr]'
$ws.Range('AE2').Value2 = 'aerobic_caustic [10^6 $/yr]'
$ws.Range('AF2').Value2 = 'TOA_fresh [10^6 $/yr]'
$ws.Range('AG2').Value2 = 'AQ336_fresh [10^6 $/yr]'
$ws.Range('AH2').Value2 = 'system_makeup_water [10^6 $/yr]'
$ws.Range('AI2').Value2 = 'Fermentation lime ratio [%]'
$ws.Range('AJ2').Value2 = 'Check [10^6 $/yr]'
$ws.Range('AK2').Value2 = 'AcrylicAcid [10^6 $/yr]'
$ws.Range('AL2').Value2 = 'ash [10^6 $/yr]'
$ws.Range('AM2').Value2 = 'gypsum [10^6 $/yr]'
$ws.Range('AN2').Value2 = 'Check [10^6 $/yr]'
$ws.Range('AO2').Value2 = 'HXN [10^9 kJ/yr]'
$ws.Range('AP2').Value2 = 'BT [10^9 kJ/yr]'
$ws.Range('AQ2').Value2 = 'Sum [10^9 kJ/yr]'
$ws.Range('AR2').Value2 = 'Check [10^9 kJ/yr]'
$ws.Range('AS2').Value2 = 'HXN [10^9 kJ/yr]'
$ws.Range('AT2').Value2 = 'CT [10^9 kJ/yr]'
$ws.Range('AU2').Value2 = 'Sum [10^9 kJ/yr]'
$ws.Range('AV2').Value2 = 'Check [10^9 kJ/yr]'
$ws.Range('AW2').Value2 = 'Sum [kW]'
$ws.Range('AX2').Value2 = 'Check [kW]'
$ws.Range('AY2').Value2 = 'Sum [10^6 $/yr]'
$ws.Range('AZ2').Value2 = 'Check [10^6 $/yr]'
$ws.Range('BA2').Value2 = 'Net present value [$]'
$ws.Range('BB2').Value2 = 'HXN energy balance error'

# Row 4
$ws.Range('A4').Value2 = 'initial'
$ws.Range("B4").Value2 = [double]"1.444711401387989"
$ws.Range("C4").Value2 = [double]"0.01839700388503976"
$ws.Range("D4").Value2 = [double]"0.9372059240906464"
$ws.Range("E4").Value2 = [double]"1.541509011224655"
$ws.Range("F4").Value2 = [double]"0.0172417810265779"
$ws.Range("G4").Value2 = [double]"0.04894712477631048"
$ws.Range("H4").Value2 = [double]"550.4327117045649"
$ws.Range("I4").Value2 = [double]"136.0671414320992"
$ws.Range("J4").Value2 = [double]"139.6488449785993"
$ws.Range("K4").Value2 = [double]"205.2578208430793"
$ws.Range("L4").Value2 = [double]"0.004346704959358991"
$ws.Range("M4").Value2 = [double]"0"
$ws.Range("N4").Value2 = [double]"132.9079607470987"
$ws.Range("O4").Value2 = [double]"527.4753828601263"
$ws.Range("P4").Value2 = [double]"0"
$ws.Range("Q4").Value2 = [double]"0"
$ws.Range("R4").Value2 = [double]"-300.05251655464"
$ws.Range("S4").Value2 = [double]"0.006547536261"
$ws.Range("T4").Value2 = [double]"0.0001877794398684067"
$ws.Range("U4").Value2 = [double]"0.003543296856379546"
$ws.Range("V4").Value2 = [double]"1.287278670316298E-06"
$ws.Range("W4").Value2 = [double]"0.001674664502029958"
$ws.Range("X4").Value2 = [double]"8.844122045762369E-06"
$ws.Range("Y4").Value2 = [double]"0.0002947832196165155"
$ws.Range("Z4").Value2 = [double]"0"
$ws.Range("AA4").Value2 = [double]"0.0002772757791182015"
$ws.Range("AB4").Value2 = [double]"1.279689809364887E-05"
$ws.Range("AC4").Value2 = [double]"0.003156830072960897"
$ws.Range("AD4").Value2 = [double]"0.000647873592065376"
$ws.Range("AE4").Value2 = [double]"0.001319431938122132"
$ws.Range("AF4").Value2 = [double]"0"
$ws.Range("AG4").Value2 = [double]"0"
$ws.Range("AH4").Value2 = [double]"3.041193949065498E-05"
$ws.Range("AI4").Value2 = [double]"0.992584721474037"
$ws.Range("AJ4").Value2 = [double]"-139.6311320352991"
$ws.Range("AK4").Value2 = [double]"0.02657836082203212"
$ws.Range("AL4").Value2 = [double]"-0.0002718148070663304"
$ws.Range("AM4").Value2 = [double]"0"
$ws.Range("AN4").Value2 = [double]"-205.2317861118714"
$ws.Range("AO4").Value2 = [double]"0"
$ws.Range("AP4").Value2 = [double]"0"
$ws.Range("AQ4").Value2 = [double]"0.4820427790017772"
$ws.Range("AR4").Value2 = [double]"0.1832166884283716"
$ws.Range("AS4").Value2 = [double]"0"
$ws.Range("AT4").Value2 = [double]"0"
$ws.Range("AU4").Value2 = [double]"-0.6629068846972399"
$ws.Range("AV4").Value2 = [double]"0"
$ws.Range("AW4").Value2 = [double]"-32770.74806031636"
$ws.Range("AX4").Value2 = [double]"32770.74806031636"
$ws.Range("AY4").Value2 = [double]"-18.0855204395274"
$ws.Range("AZ4").Value2 = [double]"18.0855204395274"
$ws.Range("BA4").Value2 = [double]"-0.1941508809104562"
$ws.Range("BB4").Value2 = [double]"-13.38891330346364"

# Row 5
$ws.Range('A5').Value2 = 'end'
$ws.Range("C5").Value2 = [double]"0.01838647561858705"
$ws.Range("D5").Value2 = [double]"0.9372059244716893"
$ws.Range("F5").Value2 = [double]"0.01723191387989405"
$ws.Range("G5").Value2 = [double]"0.04894651852624667"
$ws.Range("H5").Value2 = [double]"264.4943313167407"
$ws.Range("I5").Value2 = [double]"132.6968091273906"
$ws.Range("J5").Value2 = [double]"139.6344753134799"
$ws.Range("L5").Value2 = [double]"0.00434771713824285"
$ws.Range("M5").Value2 = [double]"0"
$ws.Range("N5").Value2 = [double]"132.9076343377146"
$ws.Range("O5").Value2 = [double]"527.3340934861478"
$ws.Range("P5").Value2 = [double]"0"
$ws.Range("Q5").Value2 = [double]"0"
$ws.Range("R5").Value2 = [double]"-138.2560409154121"
$ws.Range("S5").Value2 = [double]"0.006547536261"
$ws.Range("T5").Value2 = [double]"0.0001877794398684067"
$ws.Range("U5").Value2 = [double]"0.003543296856379546"
$ws.Range("V5").Value2 = [double]"1.287215242436646E-06"
$ws.Range("W5").Value2 = [double]"0.001674146510928734"
$ws.Range("X5").Value2 = [double]"8.84368627072649E-06"
$ws.Range("Y5").Value2 = [double]"0.0002947832196165155"
$ws.Range("Z5").Value2 = [double]"0"
$ws.Range("AA5").Value2 = [double]"0.0002772636918950845"
$ws.Range("AB5").Value2 = [double]"1.279360945898322E-05"
$ws.Range("AC5").Value2 = [double]"0.003155914920121179"
$ws.Range("AD5").Value2 = [double]"0.0006475348415891196"
$ws.Range("AE5").Value2 = [double]"0.001319389963286179"
$ws.Range("AF5").Value2 = [double]"0"
$ws.Range("AG5").Value2 = [double]"0"
$ws.Range("AH5").Value2 = [double]"3.041954664016514E-05"
$ws.Range("AI5").Value2 = [double]"0.9925824822518037"
$ws.Range("AJ5").Value2 = [double]"-139.6167641928161"
$ws.Range("AL5").Value2 = [double]"-0.000271791687012924"
$ws.Range("AM5").Value2 = [double]"0"
$ws.Range("AO5").Value2 = [double]"0"
$ws.Range("AP5").Value2 = [double]"0"
$ws.Range("AQ5").Value2 = [double]"0.482019027383671"
$ws.Range("AR5").Value2 = [double]"0.1832166884204061"
$ws.Range("AS5").Value2 = [double]"0"
$ws.Range("AT5").Value2 = [double]"0"
$ws.Range("AU5").Value2 = [double]"-0.6627365263381156"
$ws.Range("AV5").Value2 = [double]"0"
$ws.Range("AW5").Value2 = [double]"-32786.80172685236"
$ws.Range("AX5").Value2 = [double]"32786.80172685236"
$ws.Range("AY5").Value2 = [double]"-18.09438013701528"
$ws.Range("AZ5").Value2 = [double]"18.09438013701528"
$ws.Range("BB5").Value2 = [double]"-13.38295943497001"

# --- Apply the bold/bordered header style (same as the other header cells) to the
#     newly-added columns AU:BB in rows 1 and 2 ---
$newHeaderCells = "AU1","AV1","AW1","AX1","AY1","AZ1","BA1","BB1","AU2","AV2","AW2","AX2","AY2","AZ2","BA2","BB2"
foreach ($addr in $newHeaderCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

# --- Re-create the merged header cells with the new layout ---
$ws.Range("B1:Q1").Merge()
$ws.Range("S1:AJ1").Merge()
$ws.Range("AK1:AN1").Merge()
$ws.Range("AO1:AR1").Merge()
$ws.Range("AS1:AV1").Merge()
$ws.Range("AW1:AX1").Merge()
$ws.Range("AY1:AZ1").Merge()
